$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.639.37'
$ws.Range('E2').Value = '  +1.16%  '

$ws.Range('D3').Value = '3.040.41'
$ws.Range('E3').Value = '  +2.86%  '

$ws.Range('E4').Value = '  +0.18%  '

$ws.Range('D5').Value = '384.16'
$ws.Range('E5').Value = '  +1.10%  '

$ws.Range('D6').Value = '102.83'
$ws.Range('E6').Value = '  +0.63%  '

$ws.Range('D7').Value = '0.545'
$ws.Range('E7').Value = '  +0.06%  '

$ws.Range('E8').Value = '  +0.03%  '

$ws.Range('E9').Value = '  +0.13%  '

$ws.Range('D10').Value = '36.96'
$ws.Range('E10').Value = '  +0.96%  '

$ws.Range('E11').Value = '  +0.17%  '

$ws.Range('E12').Value = '  +1.11%  '

$ws.Range('D13').Value = '3.520.79'
$ws.Range('E13').Value = '  +2.90%  '

$ws.Range('D14').Value = '18.74'
$ws.Range('E14').Value = '  +2.50%  '

$ws.Range('D15').Value = '7.76'
$ws.Range('E15').Value = '  +0.16%  '

$ws.Range('D16').Value = '3.049.51'
$ws.Range('E16').Value = '  +3.15%  '

$ws.Range('D17').Value = '0.979'
$ws.Range('E17').Value = '  -1.92%  '

$ws.Range('D18').Value = '10.69'
$ws.Range('E18').Value = '  -10.40%  '

$ws.Range('D19').Value = '51.694.51'
$ws.Range('E19').Value = '  +1.19%  '

$ws.Range('E20').Value = '  +0.16%  '

$ws.Range('D21').Value = '12.44'
$ws.Range('E21').Value = '  +0.62%  '

$ws.Range('D22').Value = '0.0₃0964'
$ws.Range('E22').Value = '  +0.27%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.00'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.07%  '

$ws.Range('D24').Value = '267.61'
$ws.Range('E24').Value = '  +0.09%  '

$ws.Range('D25').Value = '3.17'
$ws.Range('E25').Value = '  -3.11%  '

$ws.Range('E26').Value = '  +5.13%  '

$ws.Range('D27').Value = '7.52'
$ws.Range('E27').Value = '  +5.06%  '

$ws.Range('E28').Value = '  +4.67%  '

$ws.Range('D29').Value = '26.37'

$ws.Range('E30').Value = '  +0.07%  '

$ws.Range('E31').Value = '  -0.89%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.30'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.85%  '

$ws.Range('E33').Value = '  +0.71%  '

$ws.Range('D34').Value = '34.14'
$ws.Range('E34').Value = '  -0.42%  '

$ws.Range('D35').Value = '50.53'
$ws.Range('E35').Value = '  -0.89%  '

$ws.Range('D36').Value = '0.0448'
$ws.Range('E36').Value = '  +3.17%  '

$ws.Range('E37').Value = '  -0.11%  '

$ws.Range('D38').Value = '3.39'
$ws.Range('E38').Value = '  +4.33%  '

$ws.Range('D39').Value = '0.287'
$ws.Range('E39').Value = '  +7.39%  '

$ws.Range('D40').Value = '17.09'
$ws.Range('E40').Value = '  +3.21%  '

$ws.Range('E41').Value = '  +1.87%  '

$ws.Range('E42').Value = '  -0.01%  '

$ws.Range('D43').Value = '127.63'
$ws.Range('E43').Value = '  +2.55%  '

$ws.Range('E44').Value = '  +1.83%  '

$ws.Range('D45').Value = '3.68'
$ws.Range('E45').Value = '  +4.12%  '

$ws.Range('D46').Value = '21.76'
$ws.Range('E46').Value = '  +1.15%  '

$ws.Range('E47').Value = '  +2.84%  '

$ws.Range('D48').Value = '2.09'
$ws.Range('E48').Value = '  +3.56%  '

$ws.Range('D49').Value = '2.039.30'
$ws.Range('E49').Value = '  -0.51%  '

$ws.Range('D50').Value = '3.342.91'
$ws.Range('E50').Value = '  +2.85%  '

$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '0.206'
$ws.Range('E51').Value = '  +6.45%  '
